# CIERRE 8 NOV 2021
# - Update the incentive-payment month text on "VALES DE INSENTIVOS" from
#   AGOSTO to OCTUBRE.
# - Move the active/selected tab from "ARQUITECTO" to "VALES DE INSENTIVOS".
# - Update the selected cell on "VALES DE INSENTIVOS" to A6.

$wb = $excel.ActiveWorkbook

$wsVales = $wb.Worksheets.Item(2)

# 1) Text content change: AGOSTO 2021 -> OCTUBRE 2021
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE OCTUBRE 2021"

# 2) Switch the active sheet/tab to "VALES DE INSENTIVOS" (this also clears
#    tabSelected on the previously-active "ARQUITECTO" sheet).
$wsVales.Activate() | Out-Null

# 3) Update the selection on the now-active sheet to A6.
$wsVales.Range("A6").Select() | Out-Null
